# Scheduled market-data refresh for Brynhildr_Profits workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per leve row
# across all class sheets, mirroring a scraped-data sync. Some rows gain/lose their
# NQ-profit (M) or HQ-profit (N) cell depending on whether that side now/no-longer
# resolves to a value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 439.85715
$ws.Range("I33").Value = 170.875
$ws.Range("J33").Value = 798.5
$ws.Range("K33").Value = 170.875
$ws.Range("L33").Value = 798.5
$ws.Range("M33").Value = 58.125
$ws.Range("N33").Value = -1256.5

$ws.Range("H64").Value = 3703.75
$ws.Range("I64").Value = 2987.4285
$ws.Range("J64").Value = 4706.6
$ws.Range("K64").Value = 2987.4285
$ws.Range("L64").Value = 4706.6
$ws.Range("M64").Value = -2739.4285
$ws.Range("N64").Value = -5202.6

$ws.Range("H67").Value = 3703.75
$ws.Range("I67").Value = 2987.4285
$ws.Range("J67").Value = 4706.6
$ws.Range("K67").Value = 2987.4285
$ws.Range("L67").Value = 4706.6
$ws.Range("M67").Value = -2129.4285
$ws.Range("N67").Value = -6422.6

$ws.Range("H107").Value = 3682.6155
$ws.Range("I107").Value = 2928.0967
$ws.Range("J107").Value = 6606.375
$ws.Range("K107").Value = 2928.0967
$ws.Range("L107").Value = 6606.375
$ws.Range("M107").Value = -1008.0967
$ws.Range("N107").Value = -10446.375

$ws.Range("H112").Value = 2882.6086
$ws.Range("J112").Value = 3000.1765
$ws.Range("L112").Value = 9000.529500000001
$ws.Range("N112").Value = -11216.5295

$ws.Range("H131").Value = 1299.6923
$ws.Range("I131").Value = 1299.6923
$ws.Range("K131").Value = 3899.0769
$ws.Range("M131").Value = 1140.9231

$ws.Range("H132").Value = 8053.7666
$ws.Range("I132").Value = 8053.7666
$ws.Range("K132").Value = 24161.2998
$ws.Range("M132").Value = -21631.2998

$ws.Range("H137").Value = 9981.883
$ws.Range("J137").Value = 12816
$ws.Range("L137").Value = 38448
$ws.Range("N137").Value = -43548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 444.46155
$ws.Range("I2").Value = 354.6
$ws.Range("J2").Value = 744
$ws.Range("K2").Value = 354.6
$ws.Range("L2").Value = 744
$ws.Range("M2").Value = -241.6
$ws.Range("N2").Value = -970

$ws.Range("H32").Value = 122506.234
$ws.Range("I32").Value = 122506.234
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 122506.234
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -122219.234
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 2187.4443
$ws.Range("J45").Value = 1999.75
$ws.Range("L45").Value = 1999.75
$ws.Range("N45").Value = -2753.75

$ws.Range("H74").Value = 11002.883
$ws.Range("I74").Value = 3184.2222
$ws.Range("J74").Value = 19798.875
$ws.Range("K74").Value = 3184.2222
$ws.Range("L74").Value = 19798.875
$ws.Range("M74").Value = -2310.2222
$ws.Range("N74").Value = -21546.875

$ws.Range("H77").Value = 11002.883
$ws.Range("I77").Value = 3184.2222
$ws.Range("J77").Value = 19798.875
$ws.Range("K77").Value = 15921.111
$ws.Range("L77").Value = 98994.375
$ws.Range("M77").Value = -11553.111
$ws.Range("N77").Value = -107730.375

$ws.Range("H116").Value = 444.46155
$ws.Range("I116").Value = 354.6
$ws.Range("J116").Value = 744
$ws.Range("K116").Value = 354.6
$ws.Range("L116").Value = 744
$ws.Range("M116").Value = 1939.4
$ws.Range("N116").Value = -5332

$ws.Range("H132").Value = 1252385.8
$ws.Range("I132").Value = 1472789.1
$ws.Range("K132").Value = 4418367.300000001
$ws.Range("M132").Value = -4415837.300000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 444.46155
$ws.Range("I3").Value = 354.6
$ws.Range("J3").Value = 744
$ws.Range("K3").Value = 354.6
$ws.Range("L3").Value = 744
$ws.Range("M3").Value = -240.6
$ws.Range("N3").Value = -972

$ws.Range("H20").Value = 4811.421
$ws.Range("I20").Value = 5338.5625
$ws.Range("K20").Value = 5338.5625
$ws.Range("M20").Value = -5091.5625

$ws.Range("H134").Value = 6543.324
$ws.Range("I134").Value = 3947.3057
$ws.Range("K134").Value = 11841.9171
$ws.Range("M134").Value = -9306.917099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4706.7915
$ws.Range("I31").Value = 5323.1763
$ws.Range("K31").Value = 5323.1763
$ws.Range("M31").Value = -5028.1763

$ws.Range("H34").Value = 4706.7915
$ws.Range("I34").Value = 5323.1763
$ws.Range("K34").Value = 5323.1763
$ws.Range("M34").Value = -5121.1763

$ws.Range("H122").Value = 8445.259
$ws.Range("I122").Value = 2423.923
$ws.Range("K122").Value = 7271.768999999999
$ws.Range("M122").Value = -4821.768999999999

$ws.Range("H132").Value = 4850.1763
$ws.Range("I132").Value = 4932.5806
$ws.Range("J132").Value = 3998.6667
$ws.Range("K132").Value = 14797.7418
$ws.Range("L132").Value = 11996.0001
$ws.Range("M132").Value = -12267.7418
$ws.Range("N132").Value = -17056.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4861.9614
$ws.Range("I129").Value = 1618.8889
$ws.Range("J129").Value = 6578.8823
$ws.Range("K129").Value = 4856.6667
$ws.Range("L129").Value = 19736.6469
$ws.Range("M129").Value = 143.3333000000002
$ws.Range("N129").Value = -29736.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 29332
$ws.Range("J63").Value = 29332
$ws.Range("L63").Value = 29332
$ws.Range("N63").Value = -30704

$ws.Range("H66").Value = 29332
$ws.Range("J66").Value = 29332
$ws.Range("L66").Value = 87996
$ws.Range("N66").Value = -94860

$ws.Range("H107").Value = 3922.8438
$ws.Range("I107").Value = 5121.909
$ws.Range("K107").Value = 5121.909
$ws.Range("M107").Value = -3201.909

$ws.Range("H132").Value = 7118.0225
$ws.Range("I132").Value = 5058.878
$ws.Range("K132").Value = 15176.634
$ws.Range("M132").Value = -12646.634

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2211.8235
$ws.Range("I22").Value = 519.6
$ws.Range("J22").Value = 2916.9167
$ws.Range("K22").Value = 519.6
$ws.Range("L22").Value = 2916.9167
$ws.Range("M22").Value = -224.6
$ws.Range("N22").Value = -3506.9167

$ws.Range("H27").Value = 2211.8235
$ws.Range("I27").Value = 519.6
$ws.Range("J27").Value = 2916.9167
$ws.Range("K27").Value = 519.6
$ws.Range("L27").Value = 2916.9167
$ws.Range("M27").Value = -412.6
$ws.Range("N27").Value = -3130.9167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15820.667
$ws.Range("I51").Value = 15820.667
$ws.Range("K51").Value = 15820.667
$ws.Range("M51").Value = -15310.667

$ws.Range("H64").Value = 49000
$ws.Range("J64").Value = 49000
$ws.Range("L64").Value = 49000
$ws.Range("N64").Value = -49496

$ws.Range("H67").Value = 49000
$ws.Range("J67").Value = 49000
$ws.Range("L67").Value = 49000
$ws.Range("N67").Value = -50716

$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 29500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 29500
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -30312

$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 29500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 29500
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -32308

$ws.Range("H96").Value = 20834936
$ws.Range("J96").Value = 2055.75
$ws.Range("L96").Value = 2055.75
$ws.Range("N96").Value = -4801.75

$ws.Range("H113").Value = 4115801.2
$ws.Range("J113").Value = 9259681
$ws.Range("L113").Value = 27779043
$ws.Range("N113").Value = -27783383
